$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.659803867340088
$ws.Range("B1").Value = 2.421384572982788
$ws.Range("C1").Value = 0.4568744599819183
$ws.Range("D1").Value = 0.4378421902656555
$ws.Range("E1").Value = 0.4637913107872009
